$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F38").Value = 443
$ws.Range("G38").Value = 16129.63
$ws.Range("F41").Value = 213
$ws.Range("G41").Value = 41085.57
$ws.Range("F42").Value = 66
$ws.Range("G42").Value = 2779.92
$ws.Range("F44").Value = 34
$ws.Range("G44").Value = 1200.88
$ws.Range("F48").Value = 238
$ws.Range("G48").Value = 13351.8
$ws.Range("F51").Value = 147
$ws.Range("G51").Value = 13750.38
$ws.Range("F53").Value = 33
$ws.Range("G53").Value = 541.53
$ws.Range("F61").Value = 234
$ws.Range("G61").Value = 61010.82
$ws.Range("B66").Value = 209237.78
$ws.Range("F144").Value = 23
$ws.Range("G144").Value = 1414.96
$ws.Range("B147").Value = 22111.99
$ws.Range("B161").Value = 57756
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 132.88
$ws.Range("F175").Value = 29
$ws.Range("G175").Value = 8410.290000000001
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("F182").Value = 21
$ws.Range("G182").Value = 1880.34
$ws.Range("F184").Value = 57
$ws.Range("G184").Value = 4674
$ws.Range("F190").Value = 7
$ws.Range("G190").Value = 623.42
$ws.Range("B193").Value = 66589.06
$ws.Range("F212").Value = 68
$ws.Range("G212").Value = 6058.12
$ws.Range("F213").Value = 219
$ws.Range("G213").Value = 27742.92
$ws.Range("F215").Value = 175
$ws.Range("G215").Value = 19650.75
$ws.Range("B218").Value = 81808.05
$ws.Range("F222").Value = 933
$ws.Range("G222").Value = 17260.5
$ws.Range("B229").Value = 29387.42
$ws.Range("F263").Value = 14
$ws.Range("G263").Value = 1451.8
$ws.Range("F264").Value = 76
$ws.Range("G264").Value = 2647.84
$ws.Range("F278").Value = 37
$ws.Range("G278").Value = 5014.98
$ws.Range("F281").Value = 3
$ws.Range("G281").Value = 80.67
$ws.Range("F288").Value = 6
$ws.Range("G288").Value = 3330.18
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F291").Value = 27
$ws.Range("G291").Value = 2313.36
$ws.Range("B292").Value = 64985
$ws.Range("C292").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F292").Value = 12
$ws.Range("G292").Value = 1052.4
$ws.Range("B293").Value = 66196
$ws.Range("C293").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F293").Value = 6
$ws.Range("G293").Value = 526.2
$ws.Range("B295").Value = 123619.18
$ws.Range("F302").Value = 37
$ws.Range("G302").Value = 4226.51
$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 78
$ws.Range("G308").Value = 3715.92
$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12
$ws.Range("B317").Value = 60325
$ws.Range("E317").Value = 151.57
$ws.Range("F317").Value = -102
$ws.Range("G317").Value = -12939.72
$ws.Range("B318").Value = 63560
$ws.Range("E318").Value = 134.87
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 126.86
$ws.Range("F324").Value = 52
$ws.Range("G324").Value = 8909.16
$ws.Range("B328").Value = -3826.03
$ws.Range("F358").Value = 46
$ws.Range("G358").Value = 10590.12
$ws.Range("F361").Value = 241
$ws.Range("G361").Value = 33882.19
$ws.Range("B363").Value = 78161.98
$ws.Range("F365").Value = 17
$ws.Range("G365").Value = 940.61
$ws.Range("F367").Value = 136
$ws.Range("G367").Value = 8255.200000000001
$ws.Range("B372").Value = 63745.66
$ws.Range("F387").Value = 436
$ws.Range("G387").Value = 42117.6
$ws.Range("B389").Value = 59059.06
$ws.Range("F408").Value = 210
$ws.Range("G408").Value = 3328.5
$ws.Range("F413").Value = 83
$ws.Range("G413").Value = 4778.31
$ws.Range("F415").Value = 58
$ws.Range("G415").Value = 3172.6
$ws.Range("B417").Value = 174015.81
$ws.Range("F429").Value = 2
$ws.Range("G429").Value = 37.56
$ws.Range("F430").Value = 230
$ws.Range("G430").Value = 10644.4
$ws.Range("F432").Value = 104
$ws.Range("G432").Value = 5034.64
$ws.Range("F433").Value = 136
$ws.Range("G433").Value = 1311.04
$ws.Range("B438").Value = 25691.02
$ws.Range("B506").Value = 64830
$ws.Range("E506").Value = 34.9
$ws.Range("F506").Value = 84
$ws.Range("G506").Value = 2757.72
$ws.Range("B507").Value = 60022
$ws.Range("E507").Value = 37.22
$ws.Range("F507").Value = -113
$ws.Range("G507").Value = -3709.79
$ws.Range("B508").Value = 41620.17
$ws.Range("F527").Value = 54
$ws.Range("G527").Value = 1787.94
$ws.Range("F529").Value = 124
$ws.Range("G529").Value = 4105.64
$ws.Range("F530").Value = 21
$ws.Range("G530").Value = 906.78
$ws.Range("F531").Value = 220
$ws.Range("G531").Value = 7284.2
$ws.Range("F532").Value = 9
$ws.Range("G532").Value = 388.62
$ws.Range("B535").Value = 24908.52
$ws.Range("F558").Value = 205
$ws.Range("G558").Value = 24979.25
$ws.Range("B561").Value = 29469.45
$ws.Range("F566").Value = 6
$ws.Range("G566").Value = 1958.76
$ws.Range("B573").Value = 27628.91
$ws.Range("F605").Value = 187
$ws.Range("G605").Value = 24889.7
$ws.Range("B607").Value = 25294.73
$ws.Range("F609").Value = 19
$ws.Range("G609").Value = 2067.39
$ws.Range("F610").Value = 64
$ws.Range("G610").Value = 1621.12
$ws.Range("F625").Value = 326
$ws.Range("G625").Value = 12006.58
$ws.Range("B628").Value = 213070.03
$ws.Range("F659").Value = 39
$ws.Range("G659").Value = 2088.06
$ws.Range("F660").Value = 51
$ws.Range("G660").Value = 1516.74
$ws.Range("F662").Value = 47
$ws.Range("G662").Value = 3774.57
$ws.Range("B668").Value = 12768.6
$ws.Range("F674").Value = 879
$ws.Range("G674").Value = 143373.69
$ws.Range("B680").Value = 144386.24
$ws.Range("F693").Value = 5
$ws.Range("G693").Value = 1211.95
$ws.Range("B713").Value = 69933.14
$ws.Range("B718").Value = 2805663.03
$ws.Range("B719").Value = 2805663.03
